$d = $word.ActiveDocument

function Replace-InRange($range, $oldText, $newText) {
    $ok = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $ok) {
        throw "Find failed for: $oldText"
    }
    return $ok
}

# --- Paragraph 6 (Objetivos PT body) ---
# "Apresentar uma abordagem pratica..." -> "Introducao. Processos fermentativos..."
$r = $d.Paragraphs.Item(6).Range
Replace-InRange $r "Apresentar uma abordagem prática da bioquímica, demonstrando as principais etapas no desenvolvimento dos processos bioquímicos industriais abordando aspectos bioquímicos importantes na produção de alimentos e importantes metabólitos. Apresentar aos alunos uma visão das aplicações potenciais e estratégicas da biotecnologia moderna, incluindo aspectos bioquímicos de bioprocessos envolvendo a utilização de hidrolisados lignocelulósicos e suas aplicações tecnológicas no contexto de biorrefinarias sustentáveis. Aprimorar o raciocínio e despertar o espírito crítico e a criatividade dos alunos na resolução de problemas industriais envolvendo processos bioquímicos." "Introdução. Processos fermentativos e enzimáticos. Processos bioquímicos industriais que incluem o processamento de alimentos, importantes metabólitos, a manufatura de bioprodutos, e os aspectos bioquímicos de bioprocessos envolvendo bioenergia e biorrefinarias."

# --- Paragraph 7 (Objetivos EN body, italic) ---
$r = $d.Paragraphs.Item(7).Range
Replace-InRange $r "Present a practical approach to biochemistry. Demonstrate the main steps in the development of industrial biochemical processes, addressing important biochemical aspects in food production, and important metabolites. Present students with a vision of the potential and strategic applications of modern biotechnology, including biochemical aspects of bioprocesses involving the use of lignocellulosic hydrolysates and their technological applications in the context of sustainable biorefineries. Improve reasoning and awaken students' critical spirit and creativity when solving industrial problems involving biochemical processes" "Introduction. Fermentative and Enzymatic processes. Industrial biochemical processes that include food processing, important metabolites, the manufacture of bioproducts, and the biochemical aspects of bioprocesses involving bioenergy and biorefineries."

# --- Paragraph 9 (Docente(s) list -> Objetivos PT body + Programa list) ---
# Replace 2nd run's text first (keeps it as a separate run), then 1st run's text.
$r = $d.Paragraphs.Item(9).Range
Replace-InRange $r "1814052 - Silvio Silverio da Silva" "1. Introdução: abordagem geral dos princípios bioquímicos aplicados em diferentes processos e setores industriais.^l2. Processos bioquímicos aplicados à indústria de alimentos: tipos de indústria de alimentos, matéria primas, fases do processamento de produtos alimentícios, conservação/alterações de alimentos 3. Processos bioquímicos nas indústrias de processamento de produtos lácteos, frutas e hortaliças, cacau e chocolate, produtos gordurosos e produtos desidratados.^l4. Principais alterações bioquímicas em alimentos, oxidação de lipídeos, escurecimento enzimático e não enzimático, controles industriais das alterações bioquímicas.^l5. Discussão e apresentação sobre aspectos bioquímicos na produção de bioprodutos de importância industrial obtidos por processos fermentativos e enzimáticos aplicáveis em diferentes setores.^l6. Bioenergia e biorrefinarias: aspectos bioquímicos de bioprocessos envolvendo a utilização de hidrolisados lignocelulósicos e suas aplicações tecnológicas"

$r = $d.Paragraphs.Item(9).Range
Replace-InRange $r "5082401 - André Moreni Lopes" "Apresentar uma abordagem prática da bioquímica, demonstrando as principais etapas no desenvolvimento dos processos bioquímicos industriais abordando aspectos bioquímicos importantes na produção de alimentos e importantes metabólitos. Apresentar aos alunos uma visão das aplicações potenciais e estratégicas da biotecnologia moderna, incluindo aspectos bioquímicos de bioprocessos envolvendo a utilização de hidrolisados lignocelulósicos e suas aplicações tecnológicas no contexto de biorrefinarias sustentáveis. Aprimorar o raciocínio e despertar o espírito crítico e a criatividade dos alunos na resolução de problemas industriais envolvendo processos bioquímicos."

# --- Paragraph 11 (Programa resumido PT body -> Avaliação Método body) ---
$r = $d.Paragraphs.Item(11).Range
Replace-InRange $r "Introdução. Processos fermentativos e enzimáticos. Processos bioquímicos industriais que incluem o processamento de alimentos, importantes metabólitos, a manufatura de bioprodutos, e os aspectos bioquímicos de bioprocessos envolvendo bioenergia e biorrefinarias." "Os alunos serão avaliados formalmente por uma prova teórica (P) e trabalhos (T). A ponderação das notas será de 70% para a prova teórica (P) e 30% para a média aritmética das notas dos trabalhos (T), ou seja: Média Final do período letivo normal (MF) = (0,7xP +0,3xT)."

# --- Paragraph 12 (Programa resumido EN body, italic -> Objetivos EN body) ---
$r = $d.Paragraphs.Item(12).Range
Replace-InRange $r "Introduction. Fermentative and Enzymatic processes. Industrial biochemical processes that include food processing, important metabolites, the manufacture of bioproducts, and the biochemical aspects of bioprocesses involving bioenergy and biorefineries." "Present a practical approach to biochemistry. Demonstrate the main steps in the development of industrial biochemical processes, addressing important biochemical aspects in food production, and important metabolites. Present students with a vision of the potential and strategic applications of modern biotechnology, including biochemical aspects of bioprocesses involving the use of lignocellulosic hydrolysates and their technological applications in the context of sustainable biorefineries. Improve reasoning and awaken students' critical spirit and creativity when solving industrial problems involving biochemical processes"

# --- Paragraph 14 (Programa PT body -> Avaliação Critério body) ---
$r = $d.Paragraphs.Item(14).Range
$oldProgramaPt = "1. Introdução: abordagem geral dos princípios bioquímicos aplicados em diferentes processos e setores industriais.^l2. Processos bioquímicos aplicados à indústria de alimentos: tipos de indústria de alimentos, matéria primas, fases do processamento de produtos alimentícios, conservação/alterações de alimentos 3. Processos bioquímicos nas indústrias de processamento de produtos lácteos, frutas e hortaliças, cacau e chocolate, produtos gordurosos e produtos desidratados.^l4. Principais alterações bioquímicas em alimentos, oxidação de lipídeos, escurecimento enzimático e não enzimático, controles industriais das alterações bioquímicas.^l5. Discussão e apresentação sobre aspectos bioquímicos na produção de bioprodutos de importância industrial obtidos por processos fermentativos e enzimáticos aplicáveis em diferentes setores.^l6. Bioenergia e biorrefinarias: aspectos bioquímicos de bioprocessos envolvendo a utilização de hidrolisados lignocelulósicos e suas aplicações tecnológicas"
Replace-InRange $r $oldProgramaPt "Serão aprovados os alunos que obtiverem média do período letivo normal igual ou maior que 5."

# --- Paragraph 17 (Avaliação: Método / Critério / Norma de recuperação bodies) ---
# Método body: the exam/works grading text -> recovery exam text
$r = $d.Paragraphs.Item(17).Range
Replace-InRange $r "Os alunos serão avaliados formalmente por uma prova teórica (P) e trabalhos (T). A ponderação das notas será de 70% para a prova teórica (P) e 30% para a média aritmética das notas dos trabalhos (T), ou seja: Média Final do período letivo normal (MF) = (0,7xP +0,3xT)." "A recuperação será feita por meio de uma prova (PR) para alunos que tenham MF maior ou igual a 3,0 e menor do que 5,0 e pelo menos 70% de frequência. A nota de recuperação (NR) será a média simples entre a média final (MF) e a prova de recuperação (PR). Será considerado aprovado o aluno com NR maior ou igual a 5,0."

# Critério body: approval-grade text -> bibliography list
$r = $d.Paragraphs.Item(17).Range
Replace-InRange $r "Serão aprovados os alunos que obtiverem média do período letivo normal igual ou maior que 5." "1.^lGAVA, A. J.; SILVA, C. A. B.; FRIAS, J. R. B. Tecnologia de alimentos - princípios e aplicações. São Paulo, Nobel, 2008. ISBN-13: 9788521313823.^l2.^lLIMA, U. A. Biotecnologia Industrial: Processos Fermentativos e Enzimáticos. Volume 3. São Paulo: Editora Edgard Blücher Ltda, 2019. ISBN 9788521214571.^l3.^lMoraes, I. O. Biotecnologia Industrial: Biotecnologia na produção de alimentos. Volume 4. São Paulo: Editora Edgard Blücher Ltda, 2021. ISBN 9786555061529"

# Norma de recuperação body: recovery exam text -> teacher name
$r = $d.Paragraphs.Item(17).Range
Replace-InRange $r "A recuperação será feita por meio de uma prova (PR) para alunos que tenham MF maior ou igual a 3,0 e menor do que 5,0 e pelo menos 70% de frequência. A nota de recuperação (NR) será a média simples entre a média final (MF) e a prova de recuperação (PR). Será considerado aprovado o aluno com NR maior ou igual a 5,0." "5082401 - André Moreni Lopes"

# --- Paragraph 19 (Bibliografia list -> second teacher name) ---
$r = $d.Paragraphs.Item(19).Range
$oldBib = "1.^lGAVA, A. J.; SILVA, C. A. B.; FRIAS, J. R. B. Tecnologia de alimentos - princípios e aplicações. São Paulo, Nobel, 2008. ISBN-13: 9788521313823.^l2.^lLIMA, U. A. Biotecnologia Industrial: Processos Fermentativos e Enzimáticos. Volume 3. São Paulo: Editora Edgard Blücher Ltda, 2019. ISBN 9788521214571.^l3.^lMoraes, I. O. Biotecnologia Industrial: Biotecnologia na produção de alimentos. Volume 4. São Paulo: Editora Edgard Blücher Ltda, 2021. ISBN 9786555061529"
Replace-InRange $r $oldBib "1814052 - Silvio Silverio da Silva"

Write-Output "Done"
